$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 184, shifting rows 184:220 down
# to 185:221 (dimension grows from A1:T220 to A1:T221).
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new weekly data point.
$ws.Range("A184").Value = 10
$ws.Range("B184").Value = "Vega Modelo de Temuco"
$ws.Range("C184").Value = "La Araucanía"
$ws.Range("D184").Value = 44798
$ws.Range("E184").Value = 9
$ws.Range("F184").Value = "Fruta"
$ws.Range("G184").Value = 100104
$ws.Range("H184").Value = "Frutos de pepita"
$ws.Range("I184").Value = 100104003
$ws.Range("J184").Value = "Membrillo"
$ws.Range("K184").Value = "Champion"
$ws.Range("L184").Value = "Primera"
$ws.Range("M184").Value = 45
$ws.Range("N184").Value = 10000
$ws.Range("O184").Value = 10000
$ws.Range("P184").Value = 10000
$ws.Range("Q184").Value = "$/bandeja 18 kilos granel"
$ws.Range("R184").Value = "Región de O'Higgins"
$ws.Range("S184").Value = 556
$ws.Range("T184").Value = 18
